# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback datetime
# strings to reflect the new report generation timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for f53f2a61-...
$wsOverview.Range("G4").Value = "2016-08-19 08:43:58"

# zh-cn sheet: Correspond Handoff / Handback DateTime for f53f2a61-...
$wsZhCn.Range("H4").Value = "2016-08-19 08:43:54"
$wsZhCn.Range("K4").Value = "2016-08-19 08:44:22"

# de-de sheet: Correspond Handoff / Handback DateTime for f53f2a61-...
$wsDeDe.Range("H4").Value = "2016-08-19 08:43:58"
$wsDeDe.Range("K4").Value = "2016-08-19 08:44:29"
